# Update "Generate Report for Handback" timestamps across sheets.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-13 17:32:44"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-13 17:32:35"
$wsZhCn.Range("K2").Value = "2016-08-13 17:33:07"

# de-de sheet: Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-13 17:33:17"
